$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H, row 1, matching the formatting
# (bold, centered, bordered) already used by the other header cells.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Column H holds the binary classification label for each patient row:
# 0 for "Control ..." rows, 1 for "MDD ..." rows (used by rows 2-11 and
# their repeated block in rows 12-21).
for ($row = 2; $row -le 21; $row++) {
    $diagnosis = $ws.Cells.Item($row, 1).Value2
    if ($diagnosis -like "MDD*") {
        $ws.Cells.Item($row, 8).Value = 1
    } else {
        $ws.Cells.Item($row, 8).Value = 0
    }
}
